$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "All studies" style Samples query for B3 (drops the Tumor / Analyte Type
# columns that the previous query selected).
$newSamplesQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001524' AND smp.sample_type = 'Blood'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# Trim the trailing newline introduced by the here-string terminator so the
# stored text ends exactly at "LIMIT 100;".
$newSamplesQuery = $newSamplesQuery.TrimEnd("`r", "`n")

$ws.Range("B3").Value = $newSamplesQuery

# Move the selection/active cell to B3 (scrolled so row 3 is visible).
$ws.Range("B3").Select()
